# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the Famfrit_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner update.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 52752.21
$ws.Range("I11").Value = 52752.21
$ws.Range("K11").Value = 52752.21
$ws.Range("M11").Value = -52612.21
$ws.Range("H132").Value = 2140.4285
$ws.Range("I132").Value = 1986.0571
$ws.Range("K132").Value = 5958.1713
$ws.Range("M132").Value = -3428.1713
$ws.Range("H137").Value = 5021
$ws.Range("I137").Value = 1914.1428
$ws.Range("J137").Value = 8127.857
$ws.Range("K137").Value = 5742.428400000001
$ws.Range("L137").Value = 24383.571
$ws.Range("M137").Value = -3192.428400000001
$ws.Range("N137").Value = -29483.571

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 8914.166999999999
$ws.Range("I8").Value = 350
$ws.Range("J8").Value = 20904
$ws.Range("K8").Value = 350
$ws.Range("L8").Value = 20904
$ws.Range("M8").Value = -206
$ws.Range("N8").Value = -21192
$ws.Range("H32").Value = 10870775
$ws.Range("I32").Value = 11495131
$ws.Range("K32").Value = 11495131
$ws.Range("M32").Value = -11494844
$ws.Range("H102").Value = 1955.4286
$ws.Range("I102").Value = 1955.4286
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1955.4286
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -333.4286
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 29420406
$ws.Range("I132").Value = 5000.3667
$ws.Range("J132").Value = 250035950
$ws.Range("K132").Value = 15001.1001
$ws.Range("L132").Value = 750107850
$ws.Range("M132").Value = -12471.1001
$ws.Range("N132").Value = -750112910

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2543.6667
$ws.Range("I99").Value = 2551.75
$ws.Range("J99").Value = 2479
$ws.Range("K99").Value = 2551.75
$ws.Range("L99").Value = 2479
$ws.Range("M99").Value = -1053.75
$ws.Range("N99").Value = -5475

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1656.6666
$ws.Range("I16").Value = 1587.1428
$ws.Range("J16").Value = 1900
$ws.Range("K16").Value = 1587.1428
$ws.Range("L16").Value = 1900
$ws.Range("M16").Value = -1300.1428
$ws.Range("N16").Value = -2474
$ws.Range("H31").Value = 35716664
$ws.Range("I31").Value = 1645.7916
$ws.Range("J31").Value = 113640340
$ws.Range("K31").Value = 1645.7916
$ws.Range("L31").Value = 113640340
$ws.Range("M31").Value = -1350.7916
$ws.Range("N31").Value = -113640930
$ws.Range("H34").Value = 35716664
$ws.Range("I34").Value = 1645.7916
$ws.Range("J34").Value = 113640340
$ws.Range("K34").Value = 1645.7916
$ws.Range("L34").Value = 113640340
$ws.Range("M34").Value = -1443.7916
$ws.Range("N34").Value = -113640744
$ws.Range("H58").Value = 8250
$ws.Range("I58").Value = 8250
$ws.Range("K58").Value = 8250
$ws.Range("M58").Value = -8047
$ws.Range("H93").Value = 20044.691
$ws.Range("I93").Value = 13612.091
$ws.Range("J93").Value = 55424
$ws.Range("K93").Value = 13612.091
$ws.Range("L93").Value = 55424
$ws.Range("M93").Value = -11740.091
$ws.Range("N93").Value = -59168
$ws.Range("H99").Value = 9611
$ws.Range("J99").Value = 9999.883
$ws.Range("L99").Value = 9999.883
$ws.Range("N99").Value = -12995.883
$ws.Range("H113").Value = 1656.6666
$ws.Range("I113").Value = 1587.1428
$ws.Range("J113").Value = 1900
$ws.Range("K113").Value = 1587.1428
$ws.Range("L113").Value = 1900
$ws.Range("M113").Value = 582.8571999999999
$ws.Range("N113").Value = -6240
$ws.Range("H126").Value = 9611
$ws.Range("J126").Value = 9999.883
$ws.Range("L126").Value = 29999.649
$ws.Range("N126").Value = -34939.649
$ws.Range("H136").Value = 8250
$ws.Range("I136").Value = 8250
$ws.Range("K136").Value = 24750
$ws.Range("M136").Value = -22200

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1267.1428
$ws.Range("I68").Value = 1500.3334
$ws.Range("J68").Value = 1203.5454
$ws.Range("K68").Value = 4501.0002
$ws.Range("L68").Value = 3610.6362
$ws.Range("M68").Value = -3690.0002
$ws.Range("N68").Value = -5232.6362
$ws.Range("H71").Value = 1267.1428
$ws.Range("I71").Value = 1500.3334
$ws.Range("J71").Value = 1203.5454
$ws.Range("K71").Value = 13503.0006
$ws.Range("L71").Value = 10831.9086
$ws.Range("M71").Value = -9447.000599999999
$ws.Range("N71").Value = -18943.9086
$ws.Range("H103").Value = 3114.25
$ws.Range("J103").Value = 3766.9092
$ws.Range("L103").Value = 11300.7276
$ws.Range("N103").Value = -13058.7276

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 15819700
$ws.Range("I11").Value = 22586714
$ws.Range("J11").Value = 29999.666
$ws.Range("K11").Value = 22586714
$ws.Range("L11").Value = 29999.666
$ws.Range("M11").Value = -22586575
$ws.Range("N11").Value = -30277.666
$ws.Range("H80").Value = 11312.167
$ws.Range("I80").Value = 10416.444
$ws.Range("J80").Value = 13999.333
$ws.Range("K80").Value = 10416.444
$ws.Range("L80").Value = 13999.333
$ws.Range("M80").Value = -9418.444
$ws.Range("N80").Value = -15995.333
$ws.Range("H83").Value = 11312.167
$ws.Range("I83").Value = 10416.444
$ws.Range("J83").Value = 13999.333
$ws.Range("K83").Value = 52082.22
$ws.Range("L83").Value = 69996.66500000001
$ws.Range("M83").Value = -47090.22
$ws.Range("N83").Value = -79980.66500000001
$ws.Range("H132").Value = 12885.436
$ws.Range("I132").Value = 11016.151
$ws.Range("K132").Value = 33048.453
$ws.Range("M132").Value = -30518.453
$ws.Range("H139").Value = 120000
$ws.Range("J139").Value = 120000
$ws.Range("L139").Value = 120000
$ws.Range("N139").Value = -130280

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1530.2703
$ws.Range("I46").Value = 987.3333
$ws.Range("K46").Value = 987.3333
$ws.Range("M46").Value = -799.3333
$ws.Range("H68").Value = 3989.9167
$ws.Range("I68").Value = 3662
$ws.Range("J68").Value = 4645.75
$ws.Range("K68").Value = 3662
$ws.Range("L68").Value = 4645.75
$ws.Range("M68").Value = -2913
$ws.Range("N68").Value = -6143.75
$ws.Range("H71").Value = 3989.9167
$ws.Range("I71").Value = 3662
$ws.Range("J71").Value = 4645.75
$ws.Range("K71").Value = 18310
$ws.Range("L71").Value = 23228.75
$ws.Range("M71").Value = -14566
$ws.Range("N71").Value = -30716.75
$ws.Range("H100").Value = 2195.2424
$ws.Range("I100").Value = 1914.9667
$ws.Range("K100").Value = 1914.9667
$ws.Range("M100").Value = -1373.9667
$ws.Range("H132").Value = 38464476
$ws.Range("I132").Value = 2972.5417
$ws.Range("J132").Value = 500002500
$ws.Range("K132").Value = 8917.625100000001
$ws.Range("L132").Value = 1500007500
$ws.Range("M132").Value = -6387.625100000001
$ws.Range("N132").Value = -1500012560
$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11870.667
$ws.Range("I62").Value = 5862.25
$ws.Range("K62").Value = 5862.25
$ws.Range("M62").Value = -5238.25
$ws.Range("H65").Value = 11870.667
$ws.Range("I65").Value = 5862.25
$ws.Range("K65").Value = 29311.25
$ws.Range("M65").Value = -26191.25
$ws.Range("H107").Value = 2386.7334
$ws.Range("I107").Value = 1925.125
$ws.Range("J107").Value = 2914.2856
$ws.Range("K107").Value = 5775.375
$ws.Range("L107").Value = 8742.856800000001
$ws.Range("M107").Value = -3855.375
$ws.Range("N107").Value = -12582.8568
$ws.Range("H113").Value = 756.5333000000001
$ws.Range("J113").Value = 799
$ws.Range("L113").Value = 2397
$ws.Range("N113").Value = -6737
$ws.Range("H126").Value = 2408.8484
$ws.Range("J126").Value = 4283.3335
$ws.Range("L126").Value = 12850.0005
$ws.Range("N126").Value = -17790.0005
$ws.Range("H132").Value = 3505.861
$ws.Range("I132").Value = 3256.625
$ws.Range("J132").Value = 5499.75
$ws.Range("K132").Value = 9769.875
$ws.Range("L132").Value = 16499.25
$ws.Range("M132").Value = -7239.875
$ws.Range("N132").Value = -21559.25
$ws.Range("H136").Value = 2177.4285
$ws.Range("I136").Value = 1707
$ws.Range("K136").Value = 5121
$ws.Range("M136").Value = -2571

